# Scheduled-runner update: refresh cached Universalis market-price /
# Leve-profit figures (columns H:N) across the ALC, ARM, BSM, CRP, GSM and
# LTW sheets. Values are plain numeric snapshots (no formulas involved).
#
# Column layout: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ,
# N=LeveProfitHQ.
#
# CRP rows 16 & 113: the HQ listing disappeared (HQ price -> 0), so the
# LeveProfitHQ figure collapses away entirely and LeveProfitNQ (column M)
# takes over the row's last value -> clear N, set M.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 63.333332
$ws.Cells.Item(4, 9).Value = 67.27273
$ws.Cells.Item(4, 10).Value = 20
$ws.Cells.Item(4, 11).Value = 67.27273
$ws.Cells.Item(4, 12).Value = 20
$ws.Cells.Item(4, 13).Value = 46.72727
$ws.Cells.Item(4, 14).Value = -248

$ws.Cells.Item(6, 8).Value = 11905073
$ws.Cells.Item(6, 9).Value = 27777922
$ws.Cells.Item(6, 10).Value = 437
$ws.Cells.Item(6, 11).Value = 83333766
$ws.Cells.Item(6, 12).Value = 1311
$ws.Cells.Item(6, 13).Value = -83333654
$ws.Cells.Item(6, 14).Value = -1535

$ws.Cells.Item(32, 8).Value = 2858.1853
$ws.Cells.Item(32, 9).Value = 3239.0667
$ws.Cells.Item(32, 10).Value = 2382.0833
$ws.Cells.Item(32, 11).Value = 3239.0667
$ws.Cells.Item(32, 12).Value = 2382.0833
$ws.Cells.Item(32, 13).Value = -2913.0667
$ws.Cells.Item(32, 14).Value = -3034.0833

$ws.Cells.Item(38, 8).Value = 2778.9412
$ws.Cells.Item(38, 9).Value = 728.5
$ws.Cells.Item(38, 11).Value = 2185.5
$ws.Cells.Item(38, 13).Value = -1813.5

$ws.Cells.Item(43, 8).Value = 5498.5
$ws.Cells.Item(43, 9).Value = 4499.5
$ws.Cells.Item(43, 10).Value = 6497.5
$ws.Cells.Item(43, 11).Value = 4499.5
$ws.Cells.Item(43, 12).Value = 6497.5
$ws.Cells.Item(43, 13).Value = -4430.5
$ws.Cells.Item(43, 14).Value = -6635.5

$ws.Cells.Item(86, 8).Value = 76959000
$ws.Cells.Item(86, 9).Value = 1766.6666
$ws.Cells.Item(86, 10).Value = 100046180
$ws.Cells.Item(86, 11).Value = 1766.6666
$ws.Cells.Item(86, 12).Value = 100046180
$ws.Cells.Item(86, 13).Value = -643.6666
$ws.Cells.Item(86, 14).Value = -100048426

$ws.Cells.Item(89, 8).Value = 76959000
$ws.Cells.Item(89, 9).Value = 1766.6666
$ws.Cells.Item(89, 10).Value = 100046180
$ws.Cells.Item(89, 11).Value = 8833.333000000001
$ws.Cells.Item(89, 12).Value = 500230900
$ws.Cells.Item(89, 13).Value = -3217.333000000001
$ws.Cells.Item(89, 14).Value = -500242132

$ws.Cells.Item(112, 8).Value = 51489.91
$ws.Cells.Item(112, 10).Value = 93361.5
$ws.Cells.Item(112, 12).Value = 280084.5
$ws.Cells.Item(112, 14).Value = -282300.5

$ws.Cells.Item(113, 8).Value = 9427.6
$ws.Cells.Item(113, 9).Value = 10899.571
$ws.Cells.Item(113, 10).Value = 5993
$ws.Cells.Item(113, 11).Value = 10899.571
$ws.Cells.Item(113, 12).Value = 5993
$ws.Cells.Item(113, 13).Value = -7645.571
$ws.Cells.Item(113, 14).Value = -12501

$ws.Cells.Item(116, 8).Value = 7414151
$ws.Cells.Item(116, 9).Value = 13894811
$ws.Cells.Item(116, 11).Value = 13894811
$ws.Cells.Item(116, 13).Value = -13891369

$ws.Cells.Item(132, 8).Value = 3593.5588
$ws.Cells.Item(132, 9).Value = 3944.423
$ws.Cells.Item(132, 10).Value = 2453.25
$ws.Cells.Item(132, 11).Value = 11833.269
$ws.Cells.Item(132, 12).Value = 7359.75
$ws.Cells.Item(132, 13).Value = -9303.269
$ws.Cells.Item(132, 14).Value = -12419.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3152.6765
$ws.Cells.Item(32, 9).Value = 3131.125
$ws.Cells.Item(32, 10).Value = 3497.5
$ws.Cells.Item(32, 11).Value = 3131.125
$ws.Cells.Item(32, 12).Value = 3497.5
$ws.Cells.Item(32, 13).Value = -2844.125
$ws.Cells.Item(32, 14).Value = -4071.5

$ws.Cells.Item(110, 8).Value = 2107.1177
$ws.Cells.Item(110, 9).Value = 1888.0667
$ws.Cells.Item(110, 11).Value = 1888.0667
$ws.Cells.Item(110, 13).Value = 156.9332999999999

$ws.Cells.Item(132, 8).Value = 3247.4119
$ws.Cells.Item(132, 9).Value = 1647.3
$ws.Cells.Item(132, 11).Value = 4941.9
$ws.Cells.Item(132, 13).Value = -2411.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 911.1177
$ws.Cells.Item(80, 9).Value = 1022.6
$ws.Cells.Item(80, 10).Value = 864.6667
$ws.Cells.Item(80, 11).Value = 1022.6
$ws.Cells.Item(80, 12).Value = 864.6667
$ws.Cells.Item(80, 13).Value = -24.60000000000002
$ws.Cells.Item(80, 14).Value = -2860.6667

$ws.Cells.Item(83, 8).Value = 911.1177
$ws.Cells.Item(83, 9).Value = 1022.6
$ws.Cells.Item(83, 10).Value = 864.6667
$ws.Cells.Item(83, 11).Value = 5113
$ws.Cells.Item(83, 12).Value = 4323.3335
$ws.Cells.Item(83, 13).Value = -121
$ws.Cells.Item(83, 14).Value = -14307.3335

$ws.Cells.Item(107, 8).Value = 3209.5715
$ws.Cells.Item(107, 10).Value = 1808
$ws.Cells.Item(107, 12).Value = 1808
$ws.Cells.Item(107, 14).Value = -5648

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 28583434
$ws.Cells.Item(6, 9).Value = 28583434
$ws.Cells.Item(6, 11).Value = 28583434
$ws.Cells.Item(6, 13).Value = -28583321

$ws.Cells.Item(16, 8).Value = 1995
$ws.Cells.Item(16, 9).Value = 1995
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1995
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -1708
$ws.Cells.Item(16, 14).ClearContents()

$ws.Cells.Item(22, 8).Value = 859.51166
$ws.Cells.Item(22, 9).Value = 850.5333000000001
$ws.Cells.Item(22, 10).Value = 880.2308
$ws.Cells.Item(22, 11).Value = 850.5333000000001
$ws.Cells.Item(22, 12).Value = 880.2308
$ws.Cells.Item(22, 13).Value = -500.5333000000001
$ws.Cells.Item(22, 14).Value = -1580.2308

$ws.Cells.Item(99, 8).Value = 6588724
$ws.Cells.Item(99, 10).Value = 4656.5713
$ws.Cells.Item(99, 12).Value = 4656.5713
$ws.Cells.Item(99, 14).Value = -7652.5713

$ws.Cells.Item(107, 8).Value = 12175.619
$ws.Cells.Item(107, 9).Value = 18505.615
$ws.Cells.Item(107, 10).Value = 1889.375
$ws.Cells.Item(107, 11).Value = 18505.615
$ws.Cells.Item(107, 12).Value = 1889.375
$ws.Cells.Item(107, 13).Value = -16585.615
$ws.Cells.Item(107, 14).Value = -5729.375

$ws.Cells.Item(113, 8).Value = 1995
$ws.Cells.Item(113, 9).Value = 1995
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 1995
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 175
$ws.Cells.Item(113, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 6588724
$ws.Cells.Item(126, 10).Value = 4656.5713
$ws.Cells.Item(126, 12).Value = 13969.7139
$ws.Cells.Item(126, 14).Value = -18909.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4289.3887
$ws.Cells.Item(70, 9).Value = 4048.0908
$ws.Cells.Item(70, 10).Value = 4668.5713
$ws.Cells.Item(70, 11).Value = 4048.0908
$ws.Cells.Item(70, 12).Value = 4668.5713
$ws.Cells.Item(70, 13).Value = -3778.0908
$ws.Cells.Item(70, 14).Value = -5208.5713

$ws.Cells.Item(73, 8).Value = 4289.3887
$ws.Cells.Item(73, 9).Value = 4048.0908
$ws.Cells.Item(73, 10).Value = 4668.5713
$ws.Cells.Item(73, 11).Value = 4048.0908
$ws.Cells.Item(73, 12).Value = 4668.5713
$ws.Cells.Item(73, 13).Value = -3112.0908
$ws.Cells.Item(73, 14).Value = -6540.5713

$ws.Cells.Item(102, 8).Value = 5345.129
$ws.Cells.Item(102, 9).Value = 6154.7915
$ws.Cells.Item(102, 11).Value = 6154.7915
$ws.Cells.Item(102, 13).Value = -4532.7915

$ws.Cells.Item(122, 8).Value = 27001.084
$ws.Cells.Item(122, 9).Value = 30716.143
$ws.Cells.Item(122, 10).Value = 21800
$ws.Cells.Item(122, 11).Value = 92148.429
$ws.Cells.Item(122, 12).Value = 65400
$ws.Cells.Item(122, 13).Value = -89698.429
$ws.Cells.Item(122, 14).Value = -70300

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 8473.352999999999
$ws.Cells.Item(93, 9).Value = 10118.538
$ws.Cells.Item(93, 10).Value = 3126.5
$ws.Cells.Item(93, 11).Value = 10118.538
$ws.Cells.Item(93, 12).Value = 3126.5
$ws.Cells.Item(93, 13).Value = -8870.538
$ws.Cells.Item(93, 14).Value = -5622.5

$ws.Cells.Item(127, 8).Value = 500300000
$ws.Cells.Item(127, 10).Value = 600000
$ws.Cells.Item(127, 12).Value = 600000
$ws.Cells.Item(127, 14).Value = -609920

